$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.707000000000001
$ws.Range("B4").Value = 6.896000000000001
$ws.Range("A11").Value = -21.464
$ws.Range("A12").Value = -21.721
$ws.Range("B14").Value = 5.938999999999999
$ws.Range("A15").Value = -21.568
$ws.Range("B26").Value = 6.706
$ws.Range("A27").Value = -20.971
$ws.Range("A28").Value = -21.163
$ws.Range("A31").Value = -21.106
$ws.Range("B31").Value = 6.205
$ws.Range("A32").Value = -20.97
$ws.Range("B35").Value = 7.273000000000001
$ws.Range("A36").Value = -20.972
$ws.Range("B37").Value = 7.320000000000002
$ws.Range("A38").Value = -19.741
$ws.Range("B39").Value = 6.883
$ws.Range("B40").Value = 8.642999999999999
$ws.Range("B45").Value = 5.89
$ws.Range("A46").Value = -21.297
$ws.Range("B52").Value = 5.46
$ws.Range("A54").Value = -21.856
$ws.Range("A55").Value = -22.214
$ws.Range("A56").Value = -21.803
$ws.Range("B57").Value = 5.332
$ws.Range("A67").Value = -21.588
$ws.Range("A69").Value = -21.721
$ws.Range("A72").Value = -21.567
$ws.Range("A73").Value = -20.628
$ws.Range("B81").Value = 6.931
$ws.Range("A83").Value = -20.146
$ws.Range("B83").Value = 6.776999999999999
$ws.Range("A86").Value = -21.911
$ws.Range("A91").Value = -21.522
$ws.Range("A93").Value = -21.49
$ws.Range("A99").Value = -20.938
$ws.Range("B100").Value = 5.517
$ws.Range("B102").Value = 6.984
